# Apply updated cryptocurrency price/volume data per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Route the literal through a formula ("="..."") and then convert the
    # cell back to a static value via copy/paste-values. This guarantees the
    # result keeps its original General/no-format style (no NumberFormat or
    # Value coercion side effects) while still landing as TEXT, matching the
    # source workbook where these columns hold plain strings (e.g. thousand-dot
    # formatted prices, "  +1.23%  " style deltas) rather than real numbers.
    $range = $ws.Range($cellRef)
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue "D2" "62.705.64"
Set-TextValue "E2" "  -2.29%  "
Set-TextValue "D3" "3.196.14"
Set-TextValue "E3" "  -3.67%  "
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "592.22"
Set-TextValue "E5" "  -1.41%  "
Set-TextValue "D6" "135.98"
Set-TextValue "E6" "  -5.58%  "
Set-TextValue "E7" "  -0.04%  "
Set-TextValue "D8" "3.201.27"
Set-TextValue "E8" "  -3.35%  "
Set-TextValue "D9" "0.506"
Set-TextValue "E9" "  -3.24%  "
Set-TextValue "D10" "0.144"
Set-TextValue "E10" "  -3.66%  "
Set-TextValue "D11" "5.34"
Set-TextValue "E11" "  -2.53%  "
Set-TextValue "D12" "0.455"
Set-TextValue "E12" "  -4.35%  "
Set-TextValue "D13" "0.0000237"
Set-TextValue "E13" "  -4.85%  "
Set-TextValue "D14" "33.50"
Set-TextValue "E14" "  -4.56%  "
Set-TextValue "D15" "3.732.18"
Set-TextValue "E15" "  -3.33%  "
Set-TextValue "E16" "  -0.05%  "
Set-TextValue "D17" "3.205.38"
Set-TextValue "E17" "  -3.14%  "
Set-TextValue "D18" "62.770.68"
Set-TextValue "E18" "  -2.31%  "
Set-TextValue "D19" "6.69"
Set-TextValue "E19" "  -3.34%  "
Set-TextValue "D20" "463.99"
Set-TextValue "E20" "  -4.50%  "
Set-TextValue "D21" "13.86"
Set-TextValue "E21" "  -3.64%  "
Set-TextValue "D22" "0.713"
Set-TextValue "E22" "  -4.39%  "
Set-TextValue "D23" "7.67"
Set-TextValue "E23" "  -5.11%  "
Set-TextValue "D24" "83.99"
Set-TextValue "E24" "  -1.24%  "
Set-TextValue "D25" "13.35"
Set-TextValue "E25" "  -2.09%  "
Set-TextValue "E26" "  -0.05%  "
Set-TextValue "D27" "2.69"
Set-TextValue "E27" "  -3.71%  "
Set-TextValue "E28" "  +0.06%  "
Set-TextValue "B29" "NEARProtocol"
Set-TextValue "C29" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D29" "6.91"
Set-TextValue "E29" "  -4.50%  "
Set-TextValue "B30" "RenderToken"
Set-TextValue "C30" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D30" "7.84"
Set-TextValue "E30" "  -5.92%  "
Set-TextValue "D31" "2.07"
Set-TextValue "E31" "  -4.64%  "
Set-TextValue "D32" "27.54"
Set-TextValue "E32" "  -3.54%  "
Set-TextValue "E33" "  -5.10%  "
Set-TextValue "D34" "2.42"
Set-TextValue "E34" "  -6.53%  "
Set-TextValue "E35" "  -4.87%  "
Set-TextValue "D36" "5.84"
Set-TextValue "E36" "  -3.24%  "
Set-TextValue "D37" "51.61"
Set-TextValue "E37" "  -3.23%  "
Set-TextValue "D38" "0.0₃0698"
Set-TextValue "E38" "  -5.52%  "
Set-TextValue "D39" "0.0390"
Set-TextValue "E39" "  -2.97%  "
Set-TextValue "D40" "419.74"
Set-TextValue "E40" "  -2.58%  "
Set-TextValue "D41" "3.005.58"
Set-TextValue "E41" "  -0.60%  "
Set-TextValue "D42" "0.115"
Set-TextValue "E42" "  +3.50%  "
Set-TextValue "D43" "8.10"
Set-TextValue "E43" "  -4.66%  "
Set-TextValue "D44" "2.61"
Set-TextValue "E44" "  -6.71%  "
Set-TextValue "D45" "0.254"
Set-TextValue "E45" "  -6.67%  "
Set-TextValue "E46" "  -5.26%  "
Set-TextValue "D48" "35.69"
Set-TextValue "E48" "  +1.72%  "
Set-TextValue "B49" "Monero"
Set-TextValue "C49" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D49" "125.10"
Set-TextValue "E49" "  +1.06%  "
Set-TextValue "B50" "InjectiveProtocol"
Set-TextValue "C50" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D50" "25.59"
Set-TextValue "E50" "  -3.22%  "
Set-TextValue "E51" "  -3.14%  "
